$p = $ppt.ActivePresentation
$s = $p.Slides.Item(9)
$sh = $s.Shapes.Item(2)
$tf = $sh.TextFrame
$tr = $tf.TextRange
$cnt = $tr.Paragraphs().Count
$para = $tr.Paragraphs($cnt)
$st = $para.Start
$len = $para.Length
$pos = $st + $len
$sub = $tr.Characters($pos, 1)
Write-Host "sub=[$($sub.Text)] start=$($sub.Start) len=$($sub.Length)"
$sub.Font.Italic = $true
Write-Host done
